$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: scattered cell edits (missing-value re-randomization) ---
# Fill in previously-missing values
$ws.Range("D3").Value = -14.2
$ws.Range("E9").Value = -6.8
$ws.Range("E10").Value = -6.1
$ws.Range("E11").Value = -7.9
$ws.Range("E12").Value = -5.3
$ws.Range("D21").Value = -14.3
# Row 33/34 (original numbering) edits -- these become rows 31/32 after the
# row deletions below, but we set them now while row numbers still match
# the source file.
$ws.Range("E33").Value = -8.1
$ws.Range("D34").Value = -14.7
$ws.Range("E34").Value = -6.4

# Clear out values that became missing
$ws.Range("E4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("D23").ClearContents()

# --- Step 2: remove the "RM 232" (row 26) and "SC 92" (row 28) records ---
# Deleting row 26 first shifts old row 28 (SC 92) up to row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()
